$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the active selection (was B3:B8 anchored at B3) to D10, as recorded
# in the saved sheet view.
$ws.Range("D10").Select()

# Rule body fix: the call site referenced the (non-existent / to-be-removed)
# brute-force "addAll" helper; point it at "addAll1" instead. The literal is
# stored as text (it starts with "=" but is not a formula), so a leading
# apostrophe is used to force text entry; Excel then marks the cell with a
# quote-prefix style so it knows to keep showing the leading apostrophe.
$ws.Range("B8").Value = "'= addAll1(null, null); `"Hello`";"

